$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Wnt5a"
$ws.Range("C2").Value = "Fzd3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1599003333333333
$ws.Range("H2").Value = 0.479701
$ws.Range("I2").Value = 0.0264777194346773
$ws.Range("J2").Value = 0.02647771943467731
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2590246666666667
$ws.Range("N2").Value = 0.777074
$ws.Range("O2").Value = 0.1073177818850196
$ws.Range("P2").Value = 0.1073177818850196
$ws.Range("Q2").Value = 0.04141813054155556
$ws.Range("R2").Value = 0.372763174874
$ws.Range("S2").Value = 0.002841530119103444
$ws.Range("T2").Value = 0.002841530119103444

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Wnt5a"
$ws.Range("C3").Value = "Fzd3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1599003333333333
$ws.Range("H3").Value = 0.479701
$ws.Range("I3").Value = 0.0264777194346773
$ws.Range("J3").Value = 0.02647771943467731
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.8886716666666666
$ws.Range("N3").Value = 2.666015
$ws.Range("O3").Value = 0.3681899230603399
$ws.Range("P3").Value = 0.3681899230603398
$ws.Range("Q3").Value = 0.1420988957238889
$ws.Range("R3").Value = 1.278890061515
$ws.Range("S3").Value = 0.009748829481467103
$ws.Range("T3").Value = 0.009748829481467103

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Wnt5a"
$ws.Range("C4").Value = "Fzd3"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1599003333333333
$ws.Range("H4").Value = 0.479701
$ws.Range("I4").Value = 0.0264777194346773
$ws.Range("J4").Value = 0.02647771943467731
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.000138
$ws.Range("N4").Value = 0.000414
$ws.Range("O4").Value = 0.0000571754578076195
$ws.Range("P4").Value = 0.0000571754578076195
$ws.Range("Q4").Value = 0.000022066246
$ws.Range("R4").Value = 0.000198596214
$ws.Range("S4").Value = 0.000001513875730379379
$ws.Range("T4").Value = 0.000001513875730379379

$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Wnt5a"
$ws.Range("C5").Value = "Fzd3"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.1599003333333333
$ws.Range("H5").Value = 0.479701
$ws.Range("I5").Value = 0.0264777194346773
$ws.Range("J5").Value = 0.02647771943467731
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.265788666666667
$ws.Range("N5").Value = 3.797366
$ws.Range("O5").Value = 0.524435119596833
$ws.Range("P5").Value = 0.524435119596833
$ws.Range("Q5").Value = 0.2024000297295556
$ws.Range("R5").Value = 1.821600267566
$ws.Range("S5").Value = 0.01388584595837638
$ws.Range("T5").Value = 0.01388584595837638

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Wnt5a"
$ws.Range("C6").Value = "Fzd3"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 5.879152
$ws.Range("H6").Value = 17.637456
$ws.Range("I6").Value = 0.9735222805653226
$ws.Range("J6").Value = 0.9735222805653228
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.2590246666666667
$ws.Range("N6").Value = 0.777074
$ws.Range("O6").Value = 0.1073177818850196
$ws.Range("P6").Value = 0.1073177818850196
$ws.Range("Q6").Value = 1.522845387082667
$ws.Range("R6").Value = 13.705608483744
$ws.Range("S6").Value = 0.1044762517659162
$ws.Range("T6").Value = 0.1044762517659162

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Wnt5a"
$ws.Range("C7").Value = "Fzd3"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 5.879152
$ws.Range("H7").Value = 17.637456
$ws.Range("I7").Value = 0.9735222805653226
$ws.Range("J7").Value = 0.9735222805653228
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.8886716666666666
$ws.Range("N7").Value = 2.666015
$ws.Range("O7").Value = 0.3681899230603399
$ws.Range("P7").Value = 0.3681899230603398
$ws.Range("Q7").Value = 5.224635806426667
$ws.Range("R7").Value = 47.02172225784
$ws.Range("S7").Value = 0.3584410935788728
$ws.Range("T7").Value = 0.3584410935788727

$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Wnt5a"
$ws.Range("C8").Value = "Fzd3"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 5.879152
$ws.Range("H8").Value = 17.637456
$ws.Range("I8").Value = 0.9735222805653226
$ws.Range("J8").Value = 0.9735222805653228
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.000138
$ws.Range("N8").Value = 0.000414
$ws.Range("O8").Value = 0.0000571754578076195
$ws.Range("P8").Value = 0.0000571754578076195
$ws.Range("Q8").Value = 0.0008113229760000001
$ws.Range("R8").Value = 0.007301906784
$ws.Range("S8").Value = 0.00005566158207724011
$ws.Range("T8").Value = 0.00005566158207724012

$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Wnt5a"
$ws.Range("C9").Value = "Fzd3"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 5.879152
$ws.Range("H9").Value = 17.637456
$ws.Range("I9").Value = 0.9735222805653226
$ws.Range("J9").Value = 0.9735222805653228
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.265788666666667
$ws.Range("N9").Value = 3.797366
$ws.Range("O9").Value = 0.524435119596833
$ws.Range("P9").Value = 0.524435119596833
$ws.Range("Q9").Value = 7.441763971210667
$ws.Range("R9").Value = 66.975875740896
$ws.Range("S9").Value = 0.5105492736384566
$ws.Range("T9").Value = 0.5105492736384566

